$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F55").Value = 6
$ws.Range("G55").Value = 1890.78
$ws.Range("B56").Value = 4980.17
$ws.Range("F83").Value = 50
$ws.Range("G83").Value = 3339.5
$ws.Range("F84").Value = 18
$ws.Range("G84").Value = 1713.78
$ws.Range("F101").Value = 67
$ws.Range("G101").Value = 2048.86
$ws.Range("F102").Value = 7
$ws.Range("G102").Value = 683.34
$ws.Range("F107").Value = 70
$ws.Range("G107").Value = 7172.2
$ws.Range("B114").Value = 294386.08
$ws.Range("F193").Value = 334
$ws.Range("G193").Value = 21643.2
$ws.Range("F197").Value = 97
$ws.Range("G197").Value = 4509.53
$ws.Range("F198").Value = 50
$ws.Range("G198").Value = 3757
$ws.Range("B200").Value = 55646.47
$ws.Range("F208").Value = 92
$ws.Range("G208").Value = 10545.96
$ws.Range("B222").Value = 62840.45
$ws.Range("B246").Value = 64973
$ws.Range("E246").Value = 35.4
$ws.Range("F246").Value = 104
$ws.Range("G246").Value = 3463.2
$ws.Range("B247").Value = 48706
$ws.Range("E247").Value = 39.8
$ws.Range("F247").Value = -144
$ws.Range("G247").Value = -4795.2
$ws.Range("F260").Value = 56
$ws.Range("G260").Value = 1795.92
$ws.Range("F272").Value = 84
$ws.Range("G272").Value = 7366.8
$ws.Range("B274").Value = 106011.33
$ws.Range("B277").Value = 61610
$ws.Range("E277").Value = 122.71
$ws.Range("F277").Value = -58
$ws.Range("G277").Value = -5957.18
$ws.Range("B278").Value = 63565
$ws.Range("E278").Value = 109.19
$ws.Range("F278").Value = 60
$ws.Range("G278").Value = 6162.6
$ws.Range("B294").Value = 57802
$ws.Range("E294").Value = 162.71
$ws.Range("F294").Value = -79
$ws.Range("G294").Value = -11334.92
$ws.Range("B295").Value = 63571
$ws.Range("F295").Value = 9
$ws.Range("G295").Value = 1291.32
$ws.Range("B296").Value = 63531
$ws.Range("E296").Value = 152.53
$ws.Range("F296").Value = 80
$ws.Range("G296").Value = 11478.4
$ws.Range("B315").Value = 63560
$ws.Range("E315").Value = 134.87
$ws.Range("F315").Value = 1
$ws.Range("G315").Value = 126.86
$ws.Range("B316").Value = 60325
$ws.Range("E316").Value = 151.57
$ws.Range("F316").Value = -102
$ws.Range("G316").Value = -12939.72
$ws.Range("F321").Value = 146
$ws.Range("G321").Value = 8570.200000000001
$ws.Range("F322").Value = 197
$ws.Range("G322").Value = 20470.27
$ws.Range("F323").Value = 36
$ws.Range("G323").Value = 4270.68
$ws.Range("F324").Value = 56
$ws.Range("G324").Value = 3311.28
$ws.Range("F328").Value = 1280
$ws.Range("G328").Value = 26918.4
$ws.Range("B339").Value = 364570.85
$ws.Range("F344").Value = 0
$ws.Range("G344").Value = 0
$ws.Range("F345").Value = 1
$ws.Range("G345").Value = 162.29
$ws.Range("B346").Value = 11128.1
$ws.Range("F370").Value = 57
$ws.Range("G370").Value = 3060.9
$ws.Range("F375").Value = 7
$ws.Range("G375").Value = 545.4400000000001
$ws.Range("F382").Value = 170
$ws.Range("G382").Value = 7311.7
$ws.Range("F384").Value = 69
$ws.Range("G384").Value = 4852.08
$ws.Range("B395").Value = 266503.02
$ws.Range("F427").Value = 123
$ws.Range("G427").Value = 4576.83
$ws.Range("B430").Value = 58013.24
$ws.Range("F434").Value = 209
$ws.Range("G434").Value = 10830.38
$ws.Range("B448").Value = 44681.5
$ws.Range("B465").Value = 53757
$ws.Range("E465").Value = 16.08
$ws.Range("F465").Value = -159
$ws.Range("G465").Value = -2138.55
$ws.Range("B466").Value = 65069
$ws.Range("E466").Value = 14.3
$ws.Range("F466").Value = 2
$ws.Range("G466").Value = 26.9
$ws.Range("F468").Value = 128
$ws.Range("G468").Value = 1683.2
$ws.Range("F478").Value = 137
$ws.Range("G478").Value = 901.46
$ws.Range("F482").Value = 456
$ws.Range("G482").Value = 3000.48
$ws.Range("B492").Value = 2832.38
$ws.Range("F495").Value = 32
$ws.Range("G495").Value = 1001.6
$ws.Range("F500").Value = 0
$ws.Range("G500").Value = 0
$ws.Range("F501").Value = 6
$ws.Range("G501").Value = 955.92
$ws.Range("B508").Value = 18262
$ws.Range("B596").Value = 60022
$ws.Range("E596").Value = 37.22
$ws.Range("F596").Value = -113
$ws.Range("G596").Value = -3709.79
$ws.Range("B597").Value = 64830
$ws.Range("E597").Value = 34.9
$ws.Range("F597").Value = 113
$ws.Range("G597").Value = 3709.79
$ws.Range("F601").Value = 30
$ws.Range("G601").Value = 1360.2
$ws.Range("B604").Value = 1533.97
$ws.Range("F698").Value = 80
$ws.Range("G698").Value = 6524.8
$ws.Range("B705").Value = 61428
$ws.Range("D705").Value = 69.16
$ws.Range("E705").Value = 73.52
$ws.Range("F705").Value = 1
$ws.Range("G705").Value = 69.16
$ws.Range("B706").Value = 63150
$ws.Range("D706").Value = 75.68000000000001
$ws.Range("E706").Value = 80.45
$ws.Range("F706").Value = 91
$ws.Range("G706").Value = 6886.88
$ws.Range("F707").Value = 154
$ws.Range("G707").Value = 3344.88
$ws.Range("F713").Value = 476
$ws.Range("G713").Value = 64264.76
$ws.Range("F715").Value = 386
$ws.Range("G715").Value = 46594.06
$ws.Range("B716").Value = 197270.88
$ws.Range("B732").Value = 65362
$ws.Range("F732").Value = 69
$ws.Range("G732").Value = 2820.03
$ws.Range("B733").Value = 65079
$ws.Range("F733").Value = 21
$ws.Range("G733").Value = 858.27
$ws.Range("F740").Value = 115
$ws.Range("G740").Value = 6553.85
$ws.Range("B743").Value = 107164.7
$ws.Range("F768").Value = 3516
$ws.Range("G768").Value = 573494.76
$ws.Range("F771").Value = 522
$ws.Range("G771").Value = 75507.3
$ws.Range("F774").Value = 244
$ws.Range("G774").Value = 31373.52
$ws.Range("B775").Value = 872702.4300000001
$ws.Range("F787").Value = 97
$ws.Range("G787").Value = 12158.95
$ws.Range("F789").Value = 53
$ws.Range("G789").Value = 2063.29
$ws.Range("F791").Value = 71
$ws.Range("G791").Value = 2806.63
$ws.Range("B792").Value = 91209.03
$ws.Range("B793").Value = 3380459.33
$ws.Range("B794").Value = 3380459.33
